# Refactor data wrangling functions:
# shift a subset of "Role Date" values in column C back by one year (365 days).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateEdits = @{
    5  = "01/02/2017"
    7  = "03/04/2017"
    8  = "08/09/2017"
    9  = "11/12/2017"
    11 = "01/03/2017"
    12 = "12/13/2017"
    15 = "02/04/2017"
    18 = "03/05/2017"
    20 = "05/06/2017"
    21 = "04/05/2017"
    23 = "07/09/2017"
    24 = "08/10/2017"
    27 = "05/04/2017"
    28 = "09/07/2017"
    31 = "04/23/2017"
    32 = "04/16/2017"
}

foreach ($row in $dateEdits.Keys) {
    $ws.Cells.Item($row, 3).Value = $dateEdits[$row]
}

# Autofit column C so its width reflects the (unchanged) date format contents.
$ws.Columns.Item(3).EntireColumn.AutoFit() | Out-Null

# Move the active selection.
$ws.Range("C22").Select() | Out-Null
